$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z10").Value = "3.11"
